# TCGenerator.xlsx edit:
#  - Insert a new "Version" sheet between "Lisez moi" and "Generator", holding a
#    changelog (Date / Auteur / Objet) with two entries, styled header + borders.
#  - Make "Version" the active sheet/tab.
#  - Add a new "inputDate" row in the TAG sheet (LEC for KW.verifyDateValue),
#    using KW.scrollAndSetText / KW.verifyDateValue.
#  - Minor selection/view tweaks on Generator and TAG sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create & position the "Version" sheet
# ---------------------------------------------------------------------------
# NOTE: Worksheets.Add() inserts right before the active sheet and any
# previously-captured sheet references can end up stale/renamed afterwards,
# so fetch "Lisez moi" only AFTER the Add() call.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Version"

$lisezMoi = $wb.Worksheets.Item("Lisez moi")
$newSheet.Move($null, $lisezMoi)

# Re-fetch by name: Move() can leave the old reference pointing at the wrong sheet.
$vs = $wb.Worksheets.Item("Version")

$vs.Range("A1").Value = "Date"
$vs.Range("B1").Value = "Auteur"
$vs.Range("C1").Value = "Objet"

$vs.Range("A2").Value = 45012
$vs.Range("B2").Value = "JM Lafarge"
$vs.Range("C2").Value = "Création"

$vs.Range("A3").Value = 45004
$vs.Range("B3").Value = "JM Lafarge"
$vs.Range("C3").Value = "Ajout inputDate"

$vs.Range("A1:C1").Interior.ThemeColor = 6
$vs.Range("A1:C1").Interior.TintAndShade = 0.6

$vs.Range("A1:C12").Borders.LineStyle = 1

$vs.Range("A2:A12").NumberFormat = "m/d/yyyy"

$vs.Columns.Item(1).ColumnWidth = 11.5546875
$vs.Columns.Item(2).ColumnWidth = 13.21875
$vs.Columns.Item(3).ColumnWidth = 78.77734375

$vs.Range("C30").Select()

# ---------------------------------------------------------------------------
# 2. Activate the "Version" tab
# ---------------------------------------------------------------------------
$vs.Activate()

# ---------------------------------------------------------------------------
# 3. TAG sheet: insert the "inputDate" LEC row (before the "inputSR" row)
# ---------------------------------------------------------------------------
$tag = $wb.Worksheets.Item("TAG")
$tag.Rows.Item(7).Insert()

$tag.Range("A7").Value = "inputDate"
$tag.Range("B7").Value = "KW.scrollAndSetText(myJDD,"""
$tag.Range("C7").Value = """)"
$tag.Range("D7").Value = "KW.scrollAndSetText(myJDD, """
$tag.Range("E7").Value = """)"
$tag.Range("F7").Value = "KW.verifyDateValue(myJDD,"""
$tag.Range("G7").Value = """)"

$tag.Range("A7:G7").Font.Size = $tag.Range("A6:G6").Font.Size
$tag.Range("A7:G7").Style = $tag.Range("D7:D7").Style

$tag.Range("A7").Select()
$tag.Range("A7:XFD7").Select()

# ---------------------------------------------------------------------------
# 4. Generator sheet: restore the previously selected (non-tab-selected) view
# ---------------------------------------------------------------------------
$gen = $wb.Worksheets.Item("Generator")
$gen.Range("B18").Select()

# ---------------------------------------------------------------------------
# 5. Leave "Version" as the active/selected tab (matches activeTab="1")
# ---------------------------------------------------------------------------
$vs = $wb.Worksheets.Item("Version")
$vs.Activate()
$vs.Range("C30").Select()
